$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set column E (Obrigatorio) to "S" for rows 2 through 11
$ws.Range("E2:E11").Value = "S"
